{"js": "// Edit 1: \"Designed and executed over 25,000 Email, Journey, and Automation test\n// cases, ensuring seamless functionality of ... UTM parameters and fallbacks\n// within each email, in accordance with the requestor's specifications.\"\n// becomes:\n// \"Designed and executed thousands of manual test cases for Emails, Journeys,\n// and Automations, ensuring seamless functionality of ... UTM parameters and\n// fallbacks within each email.\"\n\nconst body = context.document.body;\n\n// --- Part A: \"over 25,000 Email, Journey, and Automation test cases, \" ->\n//             \"thousands of manual test cases for Emails, Journeys, and Automations, \"\nconst hits1 = body.search(\"over 25,000 Email, Journey, and Automation test cases, \", { matchCase: true });\nhits1.load(\"items\");\nawait context.sync();\nif (hits1.items.length > 0) {\n  hits1.items[0].insertText(\n    \"thousands of manual test cases for Emails, Journeys, and Automations, \",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- Part B: drop the trailing qualifier clause ---\nconst hits2 = body.search(\n  \" parameters and fallbacks within each email, in accordance with the requestor's specifications.\",\n  { matchCase: true }\n);\nhits2.load(\"items\");\nawait context.sync();\nif (hits2.items.length > 0) {\n  hits2.items[0].insertText(\n    \" parameters and fallbacks within each email.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// Edit 2: \"Utilized SQL queries in Salesforce Marketing Cloud\" ->\n// \"Utilized SQL queries and test scripts in Salesforce Marketing Cloud\"\nconst hits3 = body.search(\"Utilized SQL queries in Salesforce Marketing Cloud\", { matchCase: true });\nhits3.load(\"items\");\nawait context.sync();\nif (hits3.items.length > 0) {\n  hits3.items[0].insertText(\n    \"Utilized SQL queries and test scripts in Salesforce Marketing Cloud\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Edit 1a: \"over 25,000 Email, Journey, and Automation test cases, \" ->\n#          \"thousands of manual test cases for Emails, Journeys, and Automations, \"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"over 25,000 Email, Journey, and Automation test cases, \"\n$find.Replacement.Text = \"thousands of manual test cases for Emails, Journeys, and Automations, \"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# Edit 1b: drop the trailing qualifier clause\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \" parameters and fallbacks within each email, in accordance with the requestor's specifications.\"\n$find2.Replacement.Text = \" parameters and fallbacks within each email.\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n# Edit 2: \"Utilized SQL queries in Salesforce Marketing Cloud\" ->\n#         \"Utilized SQL queries and test scripts in Salesforce Marketing Cloud\"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"Utilized SQL queries in Salesforce Marketing Cloud\"\n$find3.Replacement.Text = \"Utilized SQL queries and test scripts in Salesforce Marketing Cloud\"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2)\n"}
